$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    $new = $old -replace "^Insertion de l'entity", "Mise à jour de l'entity"
    $cell.Value = $new
}

$ws.Range("B2:B16").Value = "Update_Success"

$ws.Columns.Item(1).ColumnWidth = 51.3
